$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.95"
$ws.Range("E2").Value = "'0.05%"
$ws.Range("D3").Value = "'44.36"
$ws.Range("E3").Value = "'-0.10%"
$ws.Range("D4").Value = "'5.515"
$ws.Range("E4").Value = "'-1.53%"
$ws.Range("D5").Value = "'0.08071"
$ws.Range("E5").Value = "'-0.03%"
$ws.Range("E6").Value = "'0.74%"
$ws.Range("D7").Value = "'0.9597"
$ws.Range("E7").Value = "'0.91%"
$ws.Range("D8").Value = "'0.1129"
$ws.Range("E8").Value = "'-3.51%"
$ws.Range("D9").Value = "'0.1879"
$ws.Range("E9").Value = "'1.37%"
$ws.Range("D10").Value = "'10.22"
$ws.Range("E10").Value = "'-0.14%"
$ws.Range("D11").Value = "'0.09946"
$ws.Range("E11").Value = "'2.47%"
$ws.Range("D12").Value = "'0.04701"
$ws.Range("E12").Value = "'2.83%"
$ws.Range("D13").Value = "'0.1061"
$ws.Range("E13").Value = "'-0.58%"
$ws.Range("D14").Value = "'0.001260"
$ws.Range("E14").Value = "'-1.89%"
$ws.Range("D15").Value = "'0.04102"
$ws.Range("E15").Value = "'-2.38%"
$ws.Range("D16").Value = "'0.006130"
$ws.Range("E16").Value = "'4.04%"
$ws.Range("E17").Value = "'-0.82%"
$ws.Range("E18").Value = "'2.82%"
$ws.Range("D20").Value = "'0.3315"
$ws.Range("E20").Value = "'-4.79%"
$ws.Range("D21").Value = "'0.1395"
$ws.Range("E21").Value = "'-1.16%"
$ws.Range("E22").Value = "'2.95%"
$ws.Range("D23").Value = "'0.001313"
$ws.Range("E23").Value = "'5.33%"
$ws.Range("D24").Value = "'0.004356"
$ws.Range("E24").Value = "'0.94%"
$ws.Range("E25").Value = "'7.83%"
$ws.Range("D26").Value = "'0.0003752"
$ws.Range("E26").Value = "'-5.71%"
$ws.Range("D38").Value = "'0.02639"
$ws.Range("E38").Value = "'-0.89%"
$ws.Range("D39").Value = "'0.05626"
$ws.Range("E39").Value = "'1.36%"
$ws.Range("D40").Value = "'0.007608"
$ws.Range("E40").Value = "'0.44%"
$ws.Range("D41").Value = "'0.1402"
$ws.Range("E41").Value = "'-0.12%"
$ws.Range("D42").Value = "'0.007408"
$ws.Range("E42").Value = "'-8.17%"
$ws.Range("D43").Value = "'0.001991"
$ws.Range("E43").Value = "'-1.26%"
$ws.Range("D44").Value = "'0.008724"
$ws.Range("E44").Value = "'3.92%"
$ws.Range("D45").Value = "'0.00007120"
$ws.Range("E45").Value = "'-0.77%"
$ws.Range("E46").Value = "'0.25%"
$ws.Range("D47").Value = "'0.0005817"
$ws.Range("E47").Value = "'0.09%"
$ws.Range("D48").Value = "'0.002528"
$ws.Range("E48").Value = "'11.30%"
$ws.Range("D49").Value = "'0.003503"
$ws.Range("E49").Value = "'-12.58%"
$ws.Range("D50").Value = "'0.00002107"
$ws.Range("E50").Value = "'0.25%"
$ws.Range("D51").Value = "'0.0002006"
$ws.Range("E51").Value = "'0.25%"
